$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

$ws.Range("A22").Value = 44087
$ws.Range("A22").NumberFormat = "d-mmm"

$ws.Range("B22").Value = 0.66666666666666663
$ws.Range("B22").NumberFormat = "h:mm"

$ws.Range("C22").Value = 0.875
$ws.Range("C22").NumberFormat = "h:mm"

$ws.Range("C23").Select()
